# Version 2.0.1 solucionado error espera de base de datos
#
# Update patient record data on the emergency-room intake sheet:
#  - New patient name / record number
#  - New birth date / age / birthplace / sex
#  - Nationality corrected to feminine form
#  - Emergency-contact block cleared (name, relation, address, phone)
#  - Identification document number cleared
#  - Assistance time and urgency area updated

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Patient identification -------------------------------------------------
$ws.Range("A6").Value = "ALVARADO   CORTEZ  SOFIA  MAITE"
$ws.Range("G6").Value = "/201761944"

# --- Birth date ---------------------------------------------------------
# Writing an ISO "yyyy-mm-dd" string straight into Value would make Excel's
# automatic type detection silently convert the cell into a real date
# (changing both its type and its style). To keep this a plain text cell,
# build the text with a formula (forcing a string result), copy it as a
# value into the target cell, then remove the helper cell again.
$helper = $ws.Range("ZZ1")
$helper.Formula = '=""&"2016-07-29"'
$helper.Copy()
$ws.Range("A9").PasteSpecial(-4163)  # xlPasteValues
$ws.Columns.Item($helper.Column).Delete()

$ws.Range("D9").Value = "1 AÑO 2 MESES"
$ws.Range("E9").Value = "GUATEMALA"
$ws.Range("G9").Value = "FEMENINO"

# --- Nationality ---------------------------------------------------------
$ws.Range("E11").Value = "GUATEMALTECA"

# --- Identification document number removed ------------------------------
$ws.Range("G11").Value = ""

# --- Emergency contact block cleared --------------------------------------
$ws.Range("A13").Value = ""
$ws.Range("D13").Value = ""
$ws.Range("E13").Value = ""
$ws.Range("G13").Value = ""

# --- Assistance time / urgency area ---------------------------------------
$ws.Range("D14").Value = "Hora: 12:13:15"
$ws.Range("E14").Value = "Area de urgencia: null"
